$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column: header in H1 (copy formatting from existing header G1),
# plus data values for rows 2-3.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
